# Fix manual data: trim the stray trailing space on a handful of BDRC-ID
# ("P-number") cells in column E. Because these values become duplicates of
# already-existing values elsewhere in the sheet once trimmed (e.g. "P8213"
# already exists at E25), Excel will de-duplicate the shared-string table on
# save; for the ones that don't already exist (P3214, P4CZ16780, P3285,
# P3709) a fresh shared string is effectively created as a "resource" for
# those newer attribution entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$fixes = @{
    "E11" = "P3214"
    "E12" = "P8213"
    "E14" = "P8205"
    "E18" = "P1KG8854"
    "E28" = "P8222"
    "E31" = "P4CZ16780"
    "E35" = "P3285"
    "E38" = "P3709"
}

foreach ($addr in $fixes.Keys) {
    $ws.Range($addr).Value = $fixes[$addr]
}

# Reflect where the user ended up after editing: scrolled down with the
# active cell on the last data row.
$ws.Range("E43").Select()
$excel.ActiveWindow.ScrollRow = 15
